# Update BPSK demodulator reference sheet to include a pulse filter
# (slicer lock rate) before the slicer, and adjust the QPSK 600 max
# freq offset accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 = "BPSK 300": add slicer lock rate (column I)
$ws.Range("I4").Value = 0.87

# Row 5 = "QPSK 600": halve max freq offset (column F); the dependent
# "integral max" formula in column H recalculates automatically.
$ws.Range("F5").Value = 25

# Row 6 = "BPSK 1200": add slicer lock rate (column I)
$ws.Range("I6").Value = 0.87

# Restore the cursor/selection position as last left by the author.
$ws.Range("E11").Select()
